$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Perejil / Vega Modelo de Temuco sheet gets a new weekly price-report
# row. It is inserted at row 219 (pushing the previously-existing rows
# 219..312 down to 220..313), and the new row carries its own data.
$ws.Rows.Item(219).Insert()

$ws.Range("A219").Value2 = 10
$ws.Range("B219").Value2 = "Vega Modelo de Temuco"
$ws.Range("C219").Value2 = "La Araucanía"
$ws.Range("D219").Value2 = 44704
$ws.Range("E219").Value2 = 9
$ws.Range("F219").Value2 = 100112044
$ws.Range("G219").Value2 = "Perejil"
$ws.Range("H219").Value2 = "Sin especificar"
$ws.Range("I219").Value2 = "Primera"
$ws.Range("J219").Value2 = 30
$ws.Range("K219").Value2 = 4000
$ws.Range("L219").Value2 = 4000
$ws.Range("M219").Value2 = 4000
$ws.Range("N219").Value2 = "$/docena de atados (3 kilos)"
$ws.Range("O219").Value2 = "Provincia de Cautín"
$ws.Range("P219").Value2 = 1333
$ws.Range("Q219").Value2 = 3
$ws.Range("R219").Value2 = "Hortaliza"
